# com.tutorialninja.pagefactory -> com.tutorialninja.pageobjects
# (package rename; the accompanying workbook edit swaps the Runmode=Y demo
# rows' Email/ExpectedResult pairs on the "Data" sheet and leaves the
# "Data" tab active with E3 selected.)

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- Column C: swap the two hyperlinked e-mail cells (value + format) ---
# Capture C4's current formatting onto C3 first (before C4 is touched),
# then give C3 its new value.
$ws2.Range("C4").Copy() | Out-Null
$ws2.Range("C3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws2.Range("C3").Value = "shashvat555@gmail.com"

# C4 gets a fresh (border-less, general-format) hyperlink look, then its
# new value.
$ws2.Range("C4").Style = "Hyperlink"
$ws2.Range("C4").Value = "shashvat786@gmail.com"

# --- Column E: swap the Success/Fail cells (value + format) ---
# Capture E3's current formatting onto E4 first (before E3 is touched),
# then give E4 its new value.
$ws2.Range("E3").Copy() | Out-Null
$ws2.Range("E4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws2.Range("E4").Value = "Success"

# E3 loses all direct formatting and takes the "Fail" text.
$ws2.Range("E3").ClearFormats()
$ws2.Range("E3").Value = "Fail"

$excel.CutCopyMode = $false

# --- Make the "Data" sheet the active tab with E3 selected ---
$ws2.Activate() | Out-Null
$ws2.Range("E3").Select() | Out-Null
